$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('QUIQUIA MALLQUI, CYNTHIA ANGELLINE', '06:30', '10:15', 'QUIQUIA MALLQUI, CYNTHIA ANGELLINE', '10:15', '06:30'),
    @('VEGA RIVAS, ANDREA FERNANDA', '07:00', '10:45', 'VEGA RIVAS, ANDREA FERNANDA', '10:45', '07:00'),
    @('YOVERA ROBLES, VICTOR EDUARDO', '08:00', '11:45', 'YOVERA ROBLES, VICTOR EDUARDO', '11:45', '08:00'),
    @('MEZA PEREZ, JUAN CRISTOFER', '08:00', '11:45', 'MEZA PEREZ, JUAN CRISTOFER', '11:45', '08:00'),
    @('MARTINEZ PAZ, ROCIO ESPERANZA', '09:00', '18:00', 'CARHUARICRA ESPINOZA, FIORELLA NICOLL', '12:45', '09:00'),
    @('CARHUARICRA ESPINOZA, FIORELLA NICOLL', '09:00', '12:45', 'YANQUI BRAVO, MIRIAN LUZ', '13:00', '09:15'),
    @('YANQUI BRAVO, MIRIAN LUZ', '09:15', '13:00', 'CUSI QUISPE, ANDREA ESTEFANY', '13:30', '09:45'),
    @('CUSI QUISPE, ANDREA ESTEFANY', '09:45', '13:30', 'GOMEZ ALBINO, IDALIA GIMENA', '13:45', '10:00'),
    @('GOMEZ ALBINO, IDALIA GIMENA', '10:00', '13:45', 'MONTEZUMA DEJO, EVELYN BRUNELLA', '14:00', '10:15'),
    @('MONTEZUMA DEJO, EVELYN BRUNELLA', '10:15', '14:00', 'HUAYANAY VELASCO, ATHINA', '14:00', '10:15'),
    @('HUAYANAY VELASCO, ATHINA', '10:15', '14:00', 'MORENO CANCHANYA, ROSMERY', '14:15', '10:30'),
    @('MORENO CANCHANYA, ROSMERY', '10:30', '14:15', 'RUIZ SANTOS, CIELO CRISTHINA', '14:30', '10:45'),
    @('RUIZ SANTOS, CIELO CRISTHINA', '10:45', '14:30', 'VEGA CARDENAS, ANGELICA LOURDES', '14:45', '11:00'),
    @('VEGA CARDENAS, ANGELICA LOURDES', '11:00', '14:45', 'HUAMAN HUAMANI, ALEXIS JAVIER', '14:45', '11:00'),
    @('HUAMAN HUAMANI, ALEXIS JAVIER', '11:00', '14:45', 'Del Aguila Murayari, Darla', '14:45', '11:00'),
    @('Del Aguila Murayari, Darla', '11:00', '14:45', 'QUISPE MONDRAGÓN, JUAN ALFONSO', '15:00', '11:15'),
    @('QUISPE MONDRAGÓN, JUAN ALFONSO', '11:15', '15:00', 'RIVERA CARREÑO, DIANA DESIRÉE', '17:45', '14:00'),
    @('SOTELO GONZALES, CAMILA SOFÍA', '12:30', '21:30', 'MARTINEZ PAZ, ROCIO ESPERANZA', '18:00', '09:00'),
    @('RIVERA CARREÑO, DIANA DESIRÉE', '14:00', '17:45', 'SUAREZ JARA, YENNIFER YUSSARA', '18:15', '14:30'),
    @('SUAREZ JARA, YENNIFER YUSSARA', '14:30', '18:15', 'VILCAPOMA CHILIN, JULISSA JAZMIN', '18:15', '14:30'),
    @('VILCAPOMA CHILIN, JULISSA JAZMIN', '14:30', '18:15', 'MARTICORENA LOPEZ, DAVID CARLOS', '18:30', '14:45'),
    @('MARTICORENA LOPEZ, DAVID CARLOS', '14:45', '18:30', 'MEDINA MARCELO, NAOMI ARIADNA', '18:45', '15:00'),
    @('MEDINA MARCELO, NAOMI ARIADNA', '15:00', '18:45', 'MUÑOZ SOTOMAYOR, MIRIAN RAQUEL', '19:00', '15:15'),
    @('MUÑOZ SOTOMAYOR, MIRIAN RAQUEL', '15:15', '19:00', 'VARGAS CASTRO, LOANA VICTORIA', '20:00', '16:15'),
    @('VARGAS CASTRO, LOANA VICTORIA', '16:15', '20:00', 'TORRES RAZURI, JESUS GUSTAVO SANTIAGO', '20:15', '16:30'),
    @('TORRES RAZURI, JESUS GUSTAVO SANTIAGO', '16:30', '20:15', 'YACILA GRANDEZ, RODRIGO ANDRE', '20:15', '16:30'),
    @('YACILA GRANDEZ, RODRIGO ANDRE', '16:30', '20:15', 'CHIARA LIMA, AUGUSTO SEBASTIAN', '21:00', '17:15'),
    @('CHIARA LIMA, AUGUSTO SEBASTIAN', '17:15', '21:00', 'BONILLA SÁNCHEZ, RAÚL FERNANDO', '21:15', '17:30'),
    @('BONILLA SÁNCHEZ, RAÚL FERNANDO', '17:30', '21:15', 'ALVITE CORNEJO, ANGIE LUCERO', '21:15', '17:30'),
    @('ALVITE CORNEJO, ANGIE LUCERO', '17:30', '21:15', 'MENDOZA CRUZ, LILIANA LILIANA', '21:30', '17:45'),
    @('MENDOZA CRUZ, LILIANA LILIANA', '17:45', '21:30', 'SOTELO GONZALES, CAMILA SOFÍA', '21:30', '12:30'),
    @('AYALA TAPIA, DARCIE SOL', '18:00', '21:45', 'CAPCHA YARANGO, DAVID', '21:45', '18:00'),
    @('FLORES PAREDES, LOURDES', '18:00', '21:45', 'AYALA TAPIA, DARCIE SOL', '21:45', '18:00'),
    @('CAPCHA YARANGO, DAVID', '18:00', '21:45', 'FLORES PAREDES, LOURDES', '21:45', '18:00'),
    @('CARDENAS RICAPA, FABRIZIO ESTEBAN', '18:15', '22:00', 'CARDENAS RICAPA, FABRIZIO ESTEBAN', '22:00', '18:15'),
    @('SALAS VILLANUEVA, JAMILA DASHA', '18:15', '22:00', 'VILCHEZ CUBA, JACK ANTHONY', '22:00', '18:15'),
    @('VILCHEZ CUBA, JACK ANTHONY', '18:15', '22:00', 'SALAS VILLANUEVA, JAMILA DASHA', '22:00', '18:15'),
    @('BRENIS LÁRTIGA, SEBASTIÁN', '18:30', '22:15', 'BRENIS LÁRTIGA, SEBASTIÁN', '22:15', '18:30'),
    @('CORDOVA MONTALVO, MELANY KARINA', '19:00', '22:45', 'CORDOVA MONTALVO, MELANY KARINA', '22:45', '19:00'),
    @('INGA DELGADO, CARLOS DANIEL', '19:00', '22:45', 'INGA DELGADO, CARLOS DANIEL', '22:45', '19:00')
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

$ws.Rows(42).Delete()

Write-Host "done"